$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0.7122645628334362
$ws.Cells.Item(2, 5).Value = 0.01063976098488827
$ws.Cells.Item(3, 2).Value = 134.5163962667762
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 0.4987642214596824
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(4, 2).Value = 134.5163962667762
$ws.Cells.Item(4, 3).Value = 5
$ws.Cells.Item(4, 4).Value = 0.03398050020896767
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 2).Value = 292.2372701306018
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 0.001921716223893481
$ws.Cells.Item(5, 5).Value = 0.0001865743906692699
$ws.Cells.Item(6, 2).Value = 338.2871738311746
$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 4).Value = 0.0005914237785581779
$ws.Cells.Item(6, 5).Value = 0.0001108919584796584
$ws.Cells.Item(7, 2).Value = 411.0908011030671
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 0.002302371894805456
$ws.Cells.Item(7, 5).Value = 0.0002375463066069121
$ws.Cells.Item(8, 2).Value = 426.2
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0.1029265847196399
$ws.Cells.Item(8, 5).Value = 0.00448726254908512
$ws.Cells.Item(9, 2).Value = 442.6
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 4).Value = 0.04107328262549061
$ws.Cells.Item(9, 5).Value = 0.0009087009430418279
$ws.Cells.Item(10, 2).Value = 498.0921546689657
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0.02329854005834736
$ws.Cells.Item(10, 5).Value = 0.00194154500486228
$ws.Cells.Item(11, 2).Value = 522.1
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = 0.02914010700358869
$ws.Cells.Item(11, 5).Value = 0.0007526964908676043
$ws.Cells.Item(12, 2).Value = 627.2906398511834
$ws.Cells.Item(12, 3).Value = 5
$ws.Cells.Item(12, 4).Value = 0.01266002641348959
$ws.Cells.Item(12, 5).Value = 0.00113273920541749
$ws.Cells.Item(13, 2).Value = 639.4806312796794
$ws.Cells.Item(13, 3).Value = 5
$ws.Cells.Item(13, 4).Value = 0.009109340139086099
$ws.Cells.Item(13, 5).Value = 0.000997373007929135
$ws.Cells.Item(14, 2).Value = 665.2
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 4).Value = 0.06564432074864086
$ws.Cells.Item(14, 5).Value = 0.001153187428642613
$ws.Cells.Item(15, 2).Value = 690.8
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 4).Value = 0.01092269360327826
$ws.Cells.Item(15, 5).Value = 0.0005913560231074218
$ws.Cells.Item(16, 2).Value = 699.9918964602372
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 0.002512792308611072
$ws.Cells.Item(16, 5).Value = 0.000384902480718782
$ws.Cells.Item(17, 2).Value = 728.2049555816417
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = 0.0007773338542252403
$ws.Cells.Item(17, 5).Value = 0.000155466770845048
$ws.Cells.Item(18, 2).Value = 779.0788707099858
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 0.04320687781100344
$ws.Cells.Item(18, 5).Value = 0.003304055362017911
$ws.Cells.Item(19, 2).Value = 820.1
$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(19, 4).Value = 0.01323879839570528
$ws.Cells.Item(19, 5).Value = 0.0004928435859761914
$ws.Cells.Item(20, 2).Value = 864.8217727435604
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = 0.003647852649776206
$ws.Cells.Item(20, 5).Value = 0.0004234114682775954
$ws.Cells.Item(21, 2).Value = 997.354757246059
$ws.Cells.Item(21, 3).Value = 4
$ws.Cells.Item(21, 4).Value = 0.004268312320930733
$ws.Cells.Item(21, 5).Value = 0.0005058740528510498
$ws.Cells.Item(22, 2).Value = 1011.995701247272
$ws.Cells.Item(22, 3).Value = 3
$ws.Cells.Item(22, 4).Value = 0.004783082758831489
$ws.Cells.Item(22, 5).Value = 0.001195770689707872
$ws.Cells.Item(23, 2).Value = 1073.2
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 4).Value = 0.007268349535822038
$ws.Cells.Item(23, 5).Value = 0.0005375846583847022
$ws.Cells.Item(24, 2).Value = 1082.466210769789
$ws.Cells.Item(24, 3).Value = 5
$ws.Cells.Item(24, 4).Value = 0.02505997561845997
$ws.Cells.Item(24, 5).Value = 0.001600882963341423
$ws.Cells.Item(25, 2).Value = 1105.246920920629
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 0.005533321627902068
$ws.Cells.Item(25, 5).Value = 0.0003529978393098491
$ws.Cells.Item(26, 2).Value = 1132.989830446533
$ws.Cells.Item(26, 3).Value = 3
$ws.Cells.Item(26, 4).Value = 0.00930587291695033
$ws.Cells.Item(26, 5).Value = 0.0009377823177255106
$ws.Cells.Item(27, 2).Value = 1221.150820773939
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0.02675053742603294
$ws.Cells.Item(27, 5).Value = 0.002186342001166153
$ws.Cells.Item(28, 2).Value = 1257.056227090943
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 0.019774236949361
$ws.Cells.Item(28, 5).Value = 0.001901368937438558
$ws.Cells.Item(29, 2).Value = 1278.325416780053
$ws.Cells.Item(29, 3).Value = 2
$ws.Cells.Item(29, 4).Value = 0.01413525571574973
$ws.Cells.Item(29, 5).Value = 0.0005003630341858311
$ws.Cells.Item(30, 2).Value = 1315.93448102085
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 0.007705508831167587
$ws.Cells.Item(30, 5).Value = 0.001263198169043867
$ws.Cells.Item(31, 2).Value = 1342.277962181449
$ws.Cells.Item(31, 3).Value = 2
$ws.Cells.Item(31, 4).Value = 0.003679677657209016
$ws.Cells.Item(31, 5).Value = 0.000324677440341972
$ws.Cells.Item(32, 2).Value = 1355.9
$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 4).Value = 0.01340149393926604
$ws.Cells.Item(32, 5).Value = 0.0005089174913645331
$ws.Cells.Item(33, 2).Value = 1475.089127951248
$ws.Cells.Item(33, 3).Value = 2
$ws.Cells.Item(33, 4).Value = 0.01830332173056794
$ws.Cells.Item(33, 5).Value = 0.0006639292301277737
$ws.Cells.Item(34, 2).Value = 1485.514423856088
$ws.Cells.Item(34, 3).Value = 5
$ws.Cells.Item(34, 4).Value = 0.03166765477755253
$ws.Cells.Item(34, 5).Value = 0.003258224408520749

$ws.Rows.Item(35).Delete()
